# NIT-9012523383 "Estado de Cuenta" workbook update.
#
# The author's database refresh drops the second worker entry
# (1047405822 - HELEN MARGARITA PUELLO CASTRO, periods 2506/2507) from the
# account-statement table and updates the summary figures so the sheet now
# reflects a single remaining worker/period ("parte 1 de nuevos estado de
# cuenta").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Remove the two table rows (17-18) that list HELEN MARGARITA PUELLO CASTRO's
# overdue periods; everything below shifts up automatically.
$ws.Range("17:18").EntireRow.Delete()

# Refresh the header summary block for the remaining worker/period.
$ws.Range("E11").Value = 1211      # VALOR MORA (was 115091)
$ws.Range("C13").Value = 1         # Cant. Trabajadores (was 2)
$ws.Range("F13").Value = 1         # Cant. Periodos (was 3)

# Update the remaining detail row's "Valor Mora" amount.
$ws.Range("G16").Value = 908526    # was 0
